$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E cells stay plain text so values like "24.523.81" or
# "0.00001043" are not reinterpreted/rounded as numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.523.81"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.651.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.23%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.12"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.56"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3253"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.126"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07036"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9985"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.970"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.41"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -8.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.613"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.652.17"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001043"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -7.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06592"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9979"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.58"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.933"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.70"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -8.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.53"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.489.27"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.472"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.328"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -16.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.87"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.56"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -9.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.832.63"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.26"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.187"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.056"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.719"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -16.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08445"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.661"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.18"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -11.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.204"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -7.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.272"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02236"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -6.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06026"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -9.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2071"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -6.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.102"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9983"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5899"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -8.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.788"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.60"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -8.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5622"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -8.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.54"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.946"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -8.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06920"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.191"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.96%  "
